{"js": "// Applies the resume content edits described by the commit\n// \"Fixing typos in resume, changing some formatting\":\n//   1. \"B.S. Computer Science\" -> \"b.s. computer science\"\n//   2. \"AUGUST 2021\" (end date of the Gartner Inc. role) -> \"SEPTEMBER 2022\"\n//   3. \"... for current 7 Computer Science undergraduates ...\"\n//        -> \"... for 7 former and current Computer Science undergraduates ...\"\n//   4. \"Kotlin, , C++)\" -> \"Kotlin, C++)\"   (drop the stray extra comma)\n//   5. \"... Redis, S3, Data Dog)\" -> \"... Redis, S3)\"   (drop \", Data Dog\")\n\nasync function replaceOnce(context, searchText, replacement, options) {\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, searchOptions);\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(context, \"B.S. Computer Science\", \"b.s. computer science\");\n\nawait replaceOnce(context, \"AUGUST 2021\", \"SEPTEMBER 2022\");\n\nawait replaceOnce(\n  context,\n  \"Establish a pathway to success for current 7 Computer Science undergraduates\",\n  \"Establish a pathway to success for 7 former and current Computer Science undergraduates\"\n);\n\nawait replaceOnce(context, \"Kotlin, , C++)\", \"Kotlin, C++)\");\n\nawait replaceOnce(context, \"Redis, S3, Data Dog)\", \"Redis, S3)\");\n", "ps1": "# Applies the resume content edits described by the commit\n# \"Fixing typos in resume, changing some formatting\":\n#   1. \"B.S. Computer Science\" -> \"b.s. computer science\"\n#   2. \"AUGUST 2021\" (end date of the Gartner Inc. role) -> \"SEPTEMBER 2022\"\n#   3. \"... for current 7 Computer Science undergraduates ...\"\n#        -> \"... for 7 former and current Computer Science undergraduates ...\"\n#   4. \"Kotlin, , C++)\" -> \"Kotlin, C++)\"   (drop the stray extra comma)\n#   5. \"... Redis, S3, Data Dog)\" -> \"... Redis, S3)\"   (drop \", Data Dog\")\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-Text \"B.S. Computer Science\" \"b.s. computer science\"\n\nReplace-Text \"AUGUST 2021\" \"SEPTEMBER 2022\"\n\nReplace-Text \"Establish a pathway to success for current 7 Computer Science undergraduates\" \"Establish a pathway to success for 7 former and current Computer Science undergraduates\"\n\nReplace-Text \"Kotlin, , C++)\" \"Kotlin, C++)\"\n\nReplace-Text \"Redis, S3, Data Dog)\" \"Redis, S3)\"\n"}
